# B6-PowerPoint.pptx edit
#
# 1. Three tables (on slides 14, 15 and 16) get their table style switched
#    from the deck's custom "Table_0" style to the built-in PowerPoint
#    table style {86D330F3-A1CF-49E3-A683-0C33097B1755}.
# 2. The theme colour palette used by the slide master is repointed from
#    the "Integral" (Red Violet) palette to the "Office Theme" (Office)
#    palette.

function HexToComRgb([string]$hex) {
    # PowerPoint/VBA .RGB values are stored as a COLORREF (0xBBGGRR), i.e.
    # the reverse byte order of the "RRGGBB" hex string used in OOXML.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# --- 1. Table style swap -------------------------------------------------

$newTableStyleId = "{86D330F3-A1CF-49E3-A683-0C33097B1755}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme colour scheme swap (Integral -> Office Theme) -------------

# Order matches ThemeColorScheme.Item(1..12):
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToComRgb $officeThemeColors[$i - 1]
}
